$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new "2022-Q3" row at the
#    top of the data table (row 2), pushing the existing rows down by one
#    and renumbering the index column (A) 0..4.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Inserting a whole row shifts the existing data (rows 2-5) down to rows 3-6
# and keeps their values/styles intact automatically.
$summary.Rows.Item(2).Insert()

# The freshly inserted row 2 does not carry the same per-column styling as
# the rest of the table (Excel tries to guess a style on insert); restore it
# by copying the format from row 3 (the row that used to be row 2).
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

# Fill in the new first row: 2022-Q3
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 3
$summary.Cells.Item(2, 4).Value = 0.02

# Renumber the index column for the rest of the rows (1..4).
for ($r = 3; $r -le 6; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2) Insert a brand-new worksheet "2022-Q3" right after "总计" (i.e.
#    before the current 2nd sheet "2021-Q3"), carrying the quarterly
#    fund holdings detail table. Duplicating the sibling "2021-Q3" sheet
#    (instead of adding a blank one) keeps all of its formatting
#    (sheetPr, column/header styles, page margins, ...) intact.
# ---------------------------------------------------------------------
$sibling = $wb.Worksheets.Item(2)
$sibling.Copy($sibling)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The sibling sheet only has 2 data rows, but 2022-Q3 needs 3; add the
# missing one by copying the formatting of the last existing data row down.
$q3.Range("A3:H3").Copy()
$q3.Range("A4:H4").PasteSpecial(-4122)

# Header row.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Data rows. Columns B, D, E, F, G hold text that looks numeric in the
# source data (fund codes, percentages, ...), so force them to be stored
# as text (matching the workbook's inlineStr cells) instead of letting
# Excel auto-convert them to numbers.
$rowsData = @(
    @(0, "010447", "中邮未来成长混合A", "0.43", "91.79", "3.19", "0.0137", 8),
    @(1, "562530", "华夏中证智选1000价值稳健策略ETF", "0.54", "94.32", "0.95", "0.0051", 5),
    @(2, "010448", "中邮未来成长混合C", "0.06", "91.79", "3.19", "0.0019", 8)
)

for ($i = 0; $i -lt $rowsData.Length; $i++) {
    $r = $i + 2
    $row = $rowsData[$i]

    $q3.Cells.Item($r, 1).Value = $row[0]

    $q3.Cells.Item($r, 2).NumberFormat = "@"
    $q3.Cells.Item($r, 2).Value = $row[1]

    $q3.Cells.Item($r, 3).Value = $row[2]

    $q3.Cells.Item($r, 4).NumberFormat = "@"
    $q3.Cells.Item($r, 4).Value = $row[3]

    $q3.Cells.Item($r, 5).NumberFormat = "@"
    $q3.Cells.Item($r, 5).Value = $row[4]

    $q3.Cells.Item($r, 6).NumberFormat = "@"
    $q3.Cells.Item($r, 6).Value = $row[5]

    $q3.Cells.Item($r, 7).NumberFormat = "@"
    $q3.Cells.Item($r, 7).Value = $row[6]

    $q3.Cells.Item($r, 8).Value = $row[7]
}
